$wb = $excel.ActiveWorkbook

# ALC!row5 (Leve Item ID 5503)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 233.33333
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -185

# ALC!row11 (Leve Item ID 5533)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 452.05264
$ws.Range("I11").Value = 452.05264
$ws.Range("K11").Value = 452.05264
$ws.Range("M11").Value = -312.05264

# ALC!row15 (Leve Item ID 44146)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 517.60785
$ws.Range("I15").Value = 517.60785
$ws.Range("K15").Value = 1552.82355
$ws.Range("M15").Value = -1383.82355

# ALC!row132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3650.8572
$ws.Range("I132").Value = 3592.75
$ws.Range("K132").Value = 10778.25
$ws.Range("M132").Value = -8248.25

# ALC!row135 (Leve Item ID 44047)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1029.6111
$ws.Range("I135").Value = 829.2857
$ws.Range("J135").Value = 1730.75
$ws.Range("K135").Value = 7463.571300000001
$ws.Range("L135").Value = 15576.75
$ws.Range("M135").Value = -4928.571300000001
$ws.Range("N135").Value = -20646.75

# ALC!row138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8126.6
$ws.Range("I138").Value = 9043.706
$ws.Range("K138").Value = 27131.118
$ws.Range("M138").Value = -21991.118

# ARM!row2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2814.2666
$ws.Range("I2").Value = 2142.8333
$ws.Range("K2").Value = 2142.8333
$ws.Range("M2").Value = -2029.8333

# ARM!row32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16809.191
$ws.Range("I32").Value = 9186.147999999999
$ws.Range("K32").Value = 9186.147999999999
$ws.Range("M32").Value = -8899.147999999999

# ARM!row44 (Leve Item ID 3861)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 69999
$ws.Range("J44").Value = 69999
$ws.Range("L44").Value = 69999
$ws.Range("N44").Value = -70975

# ARM!row45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2428.8
$ws.Range("I45").Value = 1813.5385
$ws.Range("K45").Value = 1813.5385
$ws.Range("M45").Value = -1436.5385

# ARM!row63 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 7572.636
$ws.Range("I63").Value = 6974.75
$ws.Range("J63").Value = 7914.2856
$ws.Range("K63").Value = 6974.75
$ws.Range("L63").Value = 7914.2856
$ws.Range("M63").Value = -6288.75
$ws.Range("N63").Value = -9286.285599999999

# ARM!row66 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 7572.636
$ws.Range("I66").Value = 6974.75
$ws.Range("J66").Value = 7914.2856
$ws.Range("K66").Value = 34873.75
$ws.Range("L66").Value = 39571.428
$ws.Range("M66").Value = -31441.75
$ws.Range("N66").Value = -46435.428

# ARM!row74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2372
$ws.Range("I74").Value = 1007.5789
$ws.Range("K74").Value = 1007.5789
$ws.Range("M74").Value = -133.5789

# ARM!row77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2372
$ws.Range("I77").Value = 1007.5789
$ws.Range("K77").Value = 5037.8945
$ws.Range("M77").Value = -669.8945000000003

# ARM!row116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2814.2666
$ws.Range("I116").Value = 2142.8333
$ws.Range("K116").Value = 2142.8333
$ws.Range("M116").Value = 151.1667000000002

# ARM!row132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2181.1738
$ws.Range("I132").Value = 1435.2858
$ws.Range("K132").Value = 4305.857400000001
$ws.Range("M132").Value = -1775.857400000001

# BSM!row3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2814.2666
$ws.Range("I3").Value = 2142.8333
$ws.Range("K3").Value = 2142.8333
$ws.Range("M3").Value = -2028.8333

# BSM!row137 (Leve Item ID 42153)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200

# CRP!row31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4927.4585
$ws.Range("I31").Value = 3635.1333
$ws.Range("K31").Value = 3635.1333
$ws.Range("M31").Value = -3340.1333

# CRP!row34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4927.4585
$ws.Range("I34").Value = 3635.1333
$ws.Range("K34").Value = 3635.1333
$ws.Range("M34").Value = -3433.1333

# CRP!row58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3818.875
$ws.Range("I58").Value = 1515.091
$ws.Range("K58").Value = 1515.091
$ws.Range("M58").Value = -1312.091

# CRP!row136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3818.875
$ws.Range("I136").Value = 1515.091
$ws.Range("K136").Value = 4545.272999999999
$ws.Range("M136").Value = -1995.272999999999

# CUL!row37 (Leve Item ID 9516)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 90000
$ws.Range("J37").Value = 90000
$ws.Range("L37").Value = 270000
$ws.Range("N37").Value = -270224

# CUL!row52 (Leve Item ID 31902)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 843.6667
$ws.Range("J52").Value = 843.6667
$ws.Range("L52").Value = 2531.0001
$ws.Range("N52").Value = -3063.0001

# CUL!row55 (Leve Item ID 4733)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 79061
$ws.Range("I55").Value = 200798.6
$ws.Range("J55").Value = 2975
$ws.Range("K55").Value = 602395.8
$ws.Range("L55").Value = 8925
$ws.Range("M55").Value = -602218.8
$ws.Range("N55").Value = -9279

# CUL!row60 (Leve Item ID 4750)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1548
$ws.Range("J60").Value = 1599
$ws.Range("L60").Value = 4797
$ws.Range("N60").Value = -5299

# CUL!row102 (Leve Item ID 19813)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 4600

# CUL!row104 (Leve Item ID 19807)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 57116.668
$ws.Range("I104").Value = 2083.3333
$ws.Range("J104").Value = 167183.33
$ws.Range("K104").Value = 6249.999899999999
$ws.Range("L104").Value = 501549.99
$ws.Range("M104").Value = -3628.999899999999
$ws.Range("N104").Value = -506791.99

# CUL!row107 (Leve Item ID 27838)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1547.1538
$ws.Range("J107").Value = 1173.5714
$ws.Range("L107").Value = 3520.7142
$ws.Range("N107").Value = -7360.7142

# CUL!row112 (Leve Item ID 27855)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2752.1667
$ws.Range("J112").Value = 2714.2856
$ws.Range("L112").Value = 8142.8568
$ws.Range("N112").Value = -10358.8568

# CUL!row113 (Leve Item ID 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3828.6667
$ws.Range("I113").Value = 2075
$ws.Range("J113").Value = 4179.4
$ws.Range("K113").Value = 6225
$ws.Range("L113").Value = 12538.2
$ws.Range("M113").Value = -4055
$ws.Range("N113").Value = -16878.2

# CUL!row115 (Leve Item ID 27861)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 1542.3334
$ws.Range("I115").Value = 1542.3334
$ws.Range("K115").Value = 4627.0002
$ws.Range("M115").Value = -3452.0002

# CUL!row118 (Leve Item ID 27872)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1565.6666
$ws.Range("I118").Value = 1549
$ws.Range("J118").Value = 1599
$ws.Range("K118").Value = 4647
$ws.Range("L118").Value = 4797
$ws.Range("M118").Value = -3404
$ws.Range("N118").Value = -7283

# CUL!row120 (Leve Item ID 27877)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 14606.75
$ws.Range("I120").Value = 5666.3335
$ws.Range("J120").Value = 15679.6
$ws.Range("K120").Value = 16999.0005
$ws.Range("L120").Value = 47038.8
$ws.Range("M120").Value = -12161.0005
$ws.Range("N120").Value = -56714.8

# CUL!row133 (Leve Item ID 44073)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 12499.75
$ws.Range("I133").Value = 4999.5
$ws.Range("K133").Value = 14998.5
$ws.Range("M133").Value = -9938.5

# LTW!row7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4949.5
$ws.Range("J7").Value = 4949.5
$ws.Range("L7").Value = 4949.5
$ws.Range("N7").Value = -5173.5

# LTW!row126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4949.5
$ws.Range("J126").Value = 4949.5
$ws.Range("L126").Value = 14848.5
$ws.Range("N126").Value = -19788.5

# LTW!row132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5099.231
$ws.Range("I132").Value = 3465
$ws.Range("K132").Value = 10395
$ws.Range("M132").Value = -7865

# WVR!row70 (Leve Item ID 11979)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 57119
$ws.Range("I70").Value = 52797.5
$ws.Range("K70").Value = 52797.5
$ws.Range("M70").Value = -52482.5

# WVR!row73 (Leve Item ID 11979)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 57119
$ws.Range("I73").Value = 52797.5
$ws.Range("K73").Value = 52797.5
$ws.Range("M73").Value = -51705.5

# WVR!row107 (Leve Item ID 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1892.2858
$ws.Range("I107").Value = 765.3333
$ws.Range("J107").Value = 2737.5
$ws.Range("K107").Value = 2295.9999
$ws.Range("L107").Value = 8212.5
$ws.Range("M107").Value = -375.9998999999998
$ws.Range("N107").Value = -12052.5

# WVR!row126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 115977.78
$ws.Range("I126").Value = 128600
$ws.Range("K126").Value = 385800
$ws.Range("M126").Value = -383330

# WVR!row132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1739
$ws.Range("I132").Value = 733.4545000000001
$ws.Range("J132").Value = 3582.5
$ws.Range("K132").Value = 2200.3635
$ws.Range("L132").Value = 10747.5
$ws.Range("M132").Value = 329.6364999999996
$ws.Range("N132").Value = -15807.5
